# Normalize "Recorded By" (column G) entries: when a comma-separated list of
# recorders starts with "system" (case-insensitive) and has more than one
# entry, move that leading "System" entry to the end of the list (swap the
# first and last items) so "System" no longer appears first.

$ws = $excel.ActiveWorkbook.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ", "
    if ($parts.Count -gt 1 -and $parts[0].ToLower() -eq "system") {
        $first = $parts[0]
        $last = $parts[$parts.Count - 1]
        $parts[0] = $last
        $parts[$parts.Count - 1] = $first
        $cell.Value = [string]::Join(", ", $parts)
    }
}
